# Add three new archival-collection rows (MCH339-1..3) to Sheet1, below the
# existing header row, matching the body-row styling already used in the
# sheet (Calibri 10, theme text colour, no fill).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- row 2 : MCH339-1 -------------------------------------------------
$ws.Cells.Item(2, 1).Value = "MCH339-1"
$ws.Cells.Item(2, 3).Value = "WINTER STUDY MATERIAL ON HISTORY, PAPER ON SA HISTORY, PAPER ON ECONOMY AND LABOUR, SPEAK, CONSENT (WOMEN), staffrider, sa outlook"
$ws.Cells.Item(2, 5).Value = "Series"
$ws.Cells.Item(2, 6).Value = "1 Box"
$ws.Cells.Item(2, 7).Value = "LOCATION: 33H | GRAP COUNT NUMER: NONE"

# --- row 3 : MCH339-2 -------------------------------------------------
$ws.Cells.Item(3, 1).Value = "MCH339-2"
$ws.Cells.Item(3, 3).Value = "REPORT ON CTPA FOR NECC, SURPLUS PEOPLE PROJECT, AFRA NEWSLETTER, UPFRONT, END CONSCRIPTION, CRISIS NEWS, NAMDA POLICY, BLACK SASH REPORT, CHALLENGE: CHURCH PEOPLE"
$ws.Cells.Item(3, 5).Value = "Series"
$ws.Cells.Item(3, 6).Value = "1 Box"
$ws.Cells.Item(3, 7).Value = "LOCATION: 33I | GRAP COUNT NUMER: NONE"

# --- row 4 : MCH339-3 -------------------------------------------------
$ws.Cells.Item(4, 1).Value = "MCH339-3"
$ws.Cells.Item(4, 3).Value = "PERSONAL- COMBAT, TRANSKEI, CISKEI, COMMUNISM ACT, NUSAS, BANTU AREAS, PEOPLE UNDER RESTRICTION ORDERS, RACE RELATIONS, OUTLOOK, UPRIGHT, MATLASEPI, BLACK CONCIOUSNESS STUDENTS"
$ws.Cells.Item(4, 5).Value = "Series"
$ws.Cells.Item(4, 6).Value = "1 Box"
$ws.Cells.Item(4, 7).Value = "LOCATION: 33I | GRAP COUNT NUMER: NONE"

# --- styling ------------------------------------------------------------
# Column B (alternativeIdentifiers) stays untouched/blank on every new row,
# so format column A and columns C:H as two separate contiguous pieces
# (keeps B out of the written cell set) covering all three rows at once.
$colA = $ws.Range("A2:A4")
$colA.Font.ThemeColor = 1
$colA.Font.Name = "Calibri"
$colA.Font.Size = 10

$restCols = $ws.Range("C2:H4")
$restCols.Font.ThemeColor = 1
$restCols.Font.Name = "Calibri"
$restCols.Font.Size = 10

# --- view state -----------------------------------------------------------
# Keep the header frozen and leave the newly entered block selected, as in
# the saved workbook.
$ws.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $true
$ws.Range("A2:J4").Select()
